# Rename parameter labels: mlam -> megp, vlam -> vegp (header row of sheet "params")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "megp"
$ws.Range("C1").Value = "vegp"

# Reformat date axis in plot.covid: move selection to C2 and widen the
# sheet's default column width slightly.
$ws.StandardWidth = 12.19140625
$ws.Range("C2").Select()
